$d = $word.ActiveDocument

# --- Step 1 -------------------------------------------------------------
# The Google-doc hyperlink text was split across three runs
# ("...edit?usp", "=", "sharing"). Use Find/Replace to merge them back
# into a single run holding the full, correct URL text.
$d.Content.Find.Execute("edit?usp", $false, $false, $false, $false, $false, `
    $true, 1, $false, "edit?usp=sharing", 2) | Out-Null
$d.Content.Find.Execute("=sharing=sharing", $false, $false, $false, $false, $false, `
    $true, 1, $false, "=sharing", 2) | Out-Null

# --- Step 2 -------------------------------------------------------------
# Add a new "Google folder:" line with its own hyperlink, right after the
# existing "Google doc link" hyperlink paragraph and before the trailing
# blank paragraph that closes the document.
#
# We insert three placeholder paragraphs immediately in front of that
# trailing (already empty) paragraph: a throw-away character ("X"), the
# "Google folder:" label, and another throw-away character ("Y") which
# will be converted into the new hyperlink.
$countBefore = $d.Paragraphs.Count
$pLast = $d.Paragraphs($countBefore)
$insertPoint = $d.Range($pLast.Range.Start, $pLast.Range.Start)
$insertPoint.InsertBefore("X`rGoogle folder:`rY`r")

# The first newly-inserted paragraph (holding "X") should end up truly
# empty, matching the blank paragraph added in the target document.
$pEmpty = $d.Paragraphs($countBefore)
$d.Range($pEmpty.Range.Start, $pEmpty.Range.Start + 1).Delete()

# The third newly-inserted paragraph (holding "Y") becomes the new
# Google Drive folder hyperlink.
$pPlaceholder = $d.Paragraphs($countBefore + 2)
$placeholderRange = $d.Range($pPlaceholder.Range.Start, $pPlaceholder.Range.Start + 1)
$placeholderRange.Text = ""
$d.Hyperlinks.Add($placeholderRange, `
    "https://drive.google.com/drive/folders/1OQE7PiZvGxOWcV4uY46UZhPVL1rqEF6j") | Out-Null
